$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Notes")

# Update the description text (row 2)
$ws.Range("A2").Value = "Description: Urban Population (%)"

# Update the source text (row 4)
$ws.Range("A4").Value = "Source: National population and Housing census 2014: Provisional Results - Uganda Bureau of Statistics"

# Insert a new row right after the source line for the source link
$ws.Rows.Item(5).Insert()
$ws.Range("A5").Value = "Source-link: http://www.ubos.org/onlinefiles/uploads/ubos/NPHC/NPHC%202014%20PROVISIONAL%20RESULTS%20REPORT.pdf"

# Update the license note (was row 13, now row 14 after the earlier insert)
$ws.Range("A14").Value = "It is licensed under a Creative Commons Attribution 4.0 International license."

# Insert a new row after the license note for the licensing link
$ws.Rows.Item(15).Insert()
$ws.Range("A15").Value = "More information on licensing is available here: https://creativecommons.org/licenses/by/4.0/"
